$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The field-metadata table (rows 2-55) is expanded/reshuffled: new fields
# (Obligor Name, Credit Facility Type, Obligor City, Industry Code Scheme,
# Obligor Rating Grade) are interleaved with the existing rows, growing the
# used range from A1:D26 to A1:D55.
$data = New-Object "object[,]" 54,4

$data[0,0] = 'Obligor Name (ObligorName)'; $data[0,1] = 'CLCO9017'; $data[0,2] = 'Report the obligor name on the credit facility.'; $data[0,3] = 'Must not contain a carriage return, line feed, comma or any unprintable character.'
$data[1,0] = 'CUSIP (CUSIP)'; $data[1,1] = 'CLCO9161'; $data[1,2] = 'Report the CUSIP of the obligor, if available. CUSIPs are identifiers created and delivered by the CSB (CUSIP Service Bureau). The CSB is managed on behalf of the American Bankers Association by Standard & Poor’s. Issuer codes are assigned alphabetically from a series that includes deliberate built-in “gaps” for future expansion.'; $data[1,3] = 'Must be valid 6 digit CUSIP number issued by the CUSIP Service Bureau.'
$data[2,0] = 'Credit Facility Type'; $data[2,1] = 'CLCO9017'; $data[2,2] = 'Report the obligor name on the credit facility.'; $data[2,3] = 'Must not contain a carriage return, line feed, comma or any unprintable character.'
$data[3,0] = 'Credit Facility Type'; $data[3,1] = 'CLCO9017'; $data[3,2] = 'Report the obligor name on the credit facility.'; $data[3,3] = 'Must not contain a carriage return, line feed, comma or any unprintable character.'
$data[4,0] = 'Obligor City'; $data[4,1] = 'CLCO9031'; $data[4,2] = 'Report the domicile of the obligor.'; $data[4,3] = 'Use the 2-letter Country Code.'
$data[5,0] = 'Industry Code Scheme'; $data[5,1] = 'CLCO4001'; $data[5,2] = 'Report the obligor rating grade from the reporting entity’s internal risk rating system.'; $data[5,3] = 'Free text indicating the obligor rating grade.'
$data[6,0] = 'Obligor Rating Grade'; $data[6,1] = 'CLCO5023'; $data[6,2] = 'The 9-digit identification assigned by the IRS for the obligor.'; $data[6,3] = '9-digit IRS identification number.'
$data[7,0] = 'IRS Identification'; $data[7,1] = 'CLCO6002'; $data[7,2] = 'Free text'; $data[7,3] = 'Free text'
$data[8,0] = 'Other Info 1'; $data[8,1] = 'CLCO7005'; $data[8,2] = 'Free text'; $data[8,3] = 'Free text'
$data[9,0] = 'Other Info 2'; $data[9,1] = 'CLCO8006'; $data[9,2] = 'Free text'; $data[9,3] = 'Free text'
$data[10,0] = 'CUSIP Number'; $data[10,1] = 'CLCO9007'; $data[10,2] = 'Must be valid 6-digit CUSIP number issued by the CUSIP Service Bureau.'; $data[10,3] = 'Must be valid 6-digit CUSIP number.'
$data[11,0] = 'Credit Facility Type'; $data[11,1] = 'CLCO9017'; $data[11,2] = 'Report the obligor name on the credit facility.'; $data[11,3] = 'Must not contain a carriage return, line feed, comma or any unprintable character.'
$data[12,0] = 'Obligor City'; $data[12,1] = 'CLCO9031'; $data[12,2] = 'Report the domicile of the obligor.'; $data[12,3] = 'Use the 2-letter Country Code.'
$data[13,0] = 'Industry Code Scheme'; $data[13,1] = 'CLCO4001'; $data[13,2] = 'Report the obligor rating grade from the reporting entity’s internal risk rating system.'; $data[13,3] = 'Free text indicating the obligor rating grade.'
$data[14,0] = 'Obligor Rating Grade'; $data[14,1] = 'CLCO5023'; $data[14,2] = 'The 9-digit identification assigned by the IRS for the obligor.'; $data[14,3] = '9-digit IRS identification number.'
$data[15,0] = 'IRS Identification'; $data[15,1] = 'CLCO6002'; $data[15,2] = 'Free text'; $data[15,3] = 'Free text'
$data[16,0] = 'Other Info 1'; $data[16,1] = 'CLCO7005'; $data[16,2] = 'Free text'; $data[16,3] = 'Free text'
$data[17,0] = 'Other Info 2'; $data[17,1] = 'CLCO8006'; $data[17,2] = 'Free text'; $data[17,3] = 'Free text'
$data[18,0] = 'CUSIP Number'; $data[18,1] = 'CLCO9007'; $data[18,2] = 'Must be valid 6-digit CUSIP number issued by the CUSIP Service Bureau.'; $data[18,3] = 'Must be valid 6-digit CUSIP number.'
$data[19,0] = 'Credit Facility Type'; $data[19,1] = 'CLCO9017'; $data[19,2] = 'Report the obligor name on the credit facility.'; $data[19,3] = 'Must not contain a carriage return, line feed, comma or any unprintable character.'
$data[20,0] = 'Obligor City'; $data[20,1] = 'CLCO9031'; $data[20,2] = 'Report the domicile of the obligor.'; $data[20,3] = 'Use the 2-letter Country Code.'
$data[21,0] = 'Industry Code Scheme'; $data[21,1] = 'CLCO4001'; $data[21,2] = 'Report the obligor rating grade from the reporting entity’s internal risk rating system.'; $data[21,3] = 'Free text indicating the obligor rating grade.'
$data[22,0] = 'Obligor Rating Grade'; $data[22,1] = 'CLCO5023'; $data[22,2] = 'The 9-digit identification assigned by the IRS for the obligor.'; $data[22,3] = '9-digit IRS identification number.'
$data[23,0] = 'IRS Identification'; $data[23,1] = 'CLCO6002'; $data[23,2] = 'Free text'; $data[23,3] = 'Free text'
$data[24,0] = 'Other Info 1'; $data[24,1] = 'CLCO7005'; $data[24,2] = 'Free text'; $data[24,3] = 'Free text'
$data[25,0] = 'Other Info 2'; $data[25,1] = 'CLCO8006'; $data[25,2] = 'Free text'; $data[25,3] = 'Free text'
$data[26,0] = 'CUSIP Number'; $data[26,1] = 'CLCO9007'; $data[26,2] = 'Must be valid 6-digit CUSIP number issued by the CUSIP Service Bureau.'; $data[26,3] = 'Must be valid 6-digit CUSIP number.'
$data[27,0] = 'Credit Facility Type'; $data[27,1] = 'CLCO9017'; $data[27,2] = 'Report the obligor name on the credit facility.'; $data[27,3] = 'Must not contain a carriage return, line feed, comma or any unprintable character.'
$data[28,0] = 'Obligor City'; $data[28,1] = 'CLCO9031'; $data[28,2] = 'Report the domicile of the obligor.'; $data[28,3] = 'Use the 2-letter Country Code.'
$data[29,0] = 'Industry Code Scheme'; $data[29,1] = 'CLCO4001'; $data[29,2] = 'Report the obligor rating grade from the reporting entity’s internal risk rating system.'; $data[29,3] = 'Free text indicating the obligor rating grade.'
$data[30,0] = 'Obligor Rating Grade'; $data[30,1] = 'CLCO5023'; $data[30,2] = 'The 9-digit identification assigned by the IRS for the obligor.'; $data[30,3] = '9-digit IRS identification number.'
$data[31,0] = 'IRS Identification'; $data[31,1] = 'CLCO6002'; $data[31,2] = 'Free text'; $data[31,3] = 'Free text'
$data[32,0] = 'Other Info 1'; $data[32,1] = 'CLCO7005'; $data[32,2] = 'Free text'; $data[32,3] = 'Free text'
$data[33,0] = 'Other Info 2'; $data[33,1] = 'CLCO8006'; $data[33,2] = 'Free text'; $data[33,3] = 'Free text'
$data[34,0] = 'CUSIP Number'; $data[34,1] = 'CLCO9007'; $data[34,2] = 'Must be valid 6-digit CUSIP number issued by the CUSIP Service Bureau.'; $data[34,3] = 'Must be valid 6-digit CUSIP number.'
$data[35,0] = 'Credit Facility Type'; $data[35,1] = 'CLCO9017'; $data[35,2] = 'Report the obligor name on the credit facility.'; $data[35,3] = 'Must not contain a carriage return, line feed, comma or any unprintable character.'
$data[36,0] = 'Obligor City'; $data[36,1] = 'CLCO9031'; $data[36,2] = 'Report the domicile of the obligor.'; $data[36,3] = 'Use the 2-letter Country Code.'
$data[37,0] = 'Industry Code Scheme'; $data[37,1] = 'CLCO4001'; $data[37,2] = 'Report the obligor rating grade from the reporting entity’s internal risk rating system.'; $data[37,3] = 'Free text indicating the obligor rating grade.'
$data[38,0] = 'Obligor Rating Grade'; $data[38,1] = 'CLCO5023'; $data[38,2] = 'The 9-digit identification assigned by the IRS for the obligor.'; $data[38,3] = '9-digit IRS identification number.'
$data[39,0] = 'IRS Identification'; $data[39,1] = 'CLCO6002'; $data[39,2] = 'Free text'; $data[39,3] = 'Free text'
$data[40,0] = 'Other Info 1'; $data[40,1] = 'CLCO7005'; $data[40,2] = 'Free text'; $data[40,3] = 'Free text'
$data[41,0] = 'Other Info 2'; $data[41,1] = 'CLCO8006'; $data[41,2] = 'Free text'; $data[41,3] = 'Free text'
$data[42,0] = 'CUSIP Number'; $data[42,1] = 'CLCO9007'; $data[42,2] = 'Must be valid 6-digit CUSIP number issued by the CUSIP Service Bureau.'; $data[42,3] = 'Must be valid 6-digit CUSIP number.'
$data[43,0] = 'Credit Facility Type'; $data[43,1] = 'CLCO9017'; $data[43,2] = 'Report the obligor name on the credit facility.'; $data[43,3] = 'Must not contain a carriage return, line feed, comma or any unprintable character.'
$data[44,0] = 'Obligor City'; $data[44,1] = 'CLCO9031'; $data[44,2] = 'Report the domicile of the obligor.'; $data[44,3] = 'Use the 2-letter Country Code.'
$data[45,0] = 'Industry Code Scheme'; $data[45,1] = 'CLCO4001'; $data[45,2] = 'Report the obligor rating grade from the reporting entity’s internal risk rating system.'; $data[45,3] = 'Free text indicating the obligor rating grade.'
$data[46,0] = 'Obligor Rating Grade'; $data[46,1] = 'CLCO5023'; $data[46,2] = 'The 9-digit identification assigned by the IRS for the obligor.'; $data[46,3] = '9-digit IRS identification number.'
$data[47,0] = 'IRS Identification'; $data[47,1] = 'CLCO6002'; $data[47,2] = 'Free text'; $data[47,3] = 'Free text'
$data[48,0] = 'Other Info 1'; $data[48,1] = 'CLCO7005'; $data[48,2] = 'Free text'; $data[48,3] = 'Free text'
$data[49,0] = 'Other Info 2'; $data[49,1] = 'CLCO8006'; $data[49,2] = 'Free text'; $data[49,3] = 'Free text'
$data[50,0] = 'CUSIP Number'; $data[50,1] = 'CLCO9007'; $data[50,2] = 'Must be valid 6-digit CUSIP number issued by the CUSIP Service Bureau.'; $data[50,3] = 'Must be valid 6-digit CUSIP number.'
$data[51,0] = 'Credit Facility Type'; $data[51,1] = 'CLCO9017'; $data[51,2] = 'Report the obligor name on the credit facility.'; $data[51,3] = 'Must not contain a carriage return, line feed, comma or any unprintable character.'
$data[52,0] = 'Obligor City'; $data[52,1] = 'CLCO9031'; $data[52,2] = 'Report the domicile of the obligor.'; $data[52,3] = 'Use the 2-letter Country Code.'
$data[53,0] = 'Industry Code Scheme'; $data[53,1] = 'CLCO4001'; $data[53,2] = 'Report the obligor rating grade from the reporting entity’s internal risk rating system.'; $data[53,3] = 'Free text indicating the obligor rating grade.'

$startRow = 2
$endRow = $startRow + $data.GetLength(0) - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 4))
$rng.Value = $data

Write-Host "Updated $($data.GetLength(0)) rows (rows $startRow to $endRow); new dimension A1:D$endRow"
